# Apply the author's edit:
#  1) Switch the three data tables (slides 14, 15, 16) from the custom
#     "Table_0" style to the built-in table style
#     {1CBFBCB2-78D6-4A12-BCB3-6D11FAD6DA5E}.
#  2) Re-colour the deck's theme (ppt/theme/theme1.xml, used by the slide
#     master) from the "Integral" / "Red Violet" palette to the default
#     Office "Office" palette - i.e. swap the theme colour scheme.

$p = $ppt.ActivePresentation

# --- 1. Table style swap on slides 14, 15 and 16 -----------------------
$newStyleId = "{1CBFBCB2-78D6-4A12-BCB3-6D11FAD6DA5E}"
foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newStyleId)
        }
    }
}

# --- 2. Theme colour scheme swap (Integral -> Office Theme colours) ----
$master = $p.SlideMaster
$theme = $master.Theme
$colorScheme = $theme.ThemeColorScheme

# Index order in ThemeColorScheme.Colors matches the OOXML <a:clrScheme>
# child order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
# RGB is passed as a standard VBA-style BGR-packed long (0xBBGGRR).
$officeColors = @(
    0x000000,  # dk1      000000
    0xFFFFFF,  # lt1      FFFFFF
    0x6A5444,  # dk2      44546A
    0xE6E6E7,  # lt2      E7E6E6
    0xD59B5B,  # accent1  5B9BD5
    0x317DED,  # accent2  ED7D31
    0xA5A5A5,  # accent3  A5A5A5
    0x00C0FF,  # accent4  FFC000
    0xC47244,  # accent5  4472C4
    0x47AD70,  # accent6  70AD47
    0xC16305,  # hlink    0563C1
    0x724F95   # folHlink 954F72
)

for ($i = 1; $i -le 12; $i++) {
    $colorScheme.Colors($i).RGB = $officeColors[$i - 1]
}
